$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 3 updates
    $ws.Range("F3").Value = 18
    $ws.Range("G3").Value = 29.9

    # Row 5 updates
    $ws.Range("D5").Value = "鼎湖路与永宁路交叉口南120米 缙云万地广场"
    $ws.Range("G5").Value = 29.9
    $ws.Range("I5").Value = "//i2.hdslb.com/bfs/openplatform/202409/vEHR9otg1726824051090.jpeg"
}
